$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Child")

$values = @{
    2  = "-2,9"
    3  = "-7,-3"
    4  = "-8,-2"
    5  = "-4,5"
    6  = "6,6"
    7  = "0,8"
    8  = "-10,9"
    9  = "3,6"
    10 = "-7,-9"
    11 = "-8,4"
    12 = "3,-1"
    13 = "3,7"
    14 = "-9,0"
    15 = "6,3"
    16 = "4,-8"
    17 = "4,7"
    18 = "-5,-5"
    19 = "-7,-1"
    20 = "-3,-4"
    21 = "5,-7"
}

foreach ($row in $values.Keys) {
    $ws.Range("D$row").Value = $values[$row]
}
